# Mark the Sets-group commands SADD, SCARD, SDIFF, SDIFFSTORE, SINTER,
# SINTERSTORE, SISMEMBER, SMEMBERS, SUNION and SUNIONSTORE as Finished,
# and record their implementing method name(s) (commit: "SMEMBERS,
# SDIFFSTORE, SINTERSTORE, SUNIONSTORE").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (Finished, StringMethod, Method)
# Order below is the original author's entry order (inferred from the
# shared-string table ordering in the target diff): rows 161-166 top to
# bottom, then 173-174, then back up to 167-168.
$rows = @(
    @{ Row = 161; StringMethod = "SetAddMemberStringAsync";               Method = "SetAddMemberAsync" }               # SADD
    @{ Row = 162; StringMethod = $null;                                   Method = "SetCardinalityAsync" }             # SCARD
    @{ Row = 163; StringMethod = "SetGetDifferenceMembersStringAsync";    Method = "SetGetDifferenceMembersAsync" }    # SDIFF
    @{ Row = 164; StringMethod = $null;                                   Method = "SetStoreDifferenceMembersAsync" } # SDIFFSTORE
    @{ Row = 165; StringMethod = "SetGetIntersectionMembersStringAsync";  Method = "SetGetIntersectionMembersAsync" }  # SINTER
    @{ Row = 166; StringMethod = $null;                                   Method = "SetStoreIntersectionMembersAsync" } # SINTERSTORE
    @{ Row = 173; StringMethod = "SetGetUnionMembersStringAsync";         Method = "SetGetUnionMembersAsync" }         # SUNION
    @{ Row = 174; StringMethod = $null;                                   Method = "SetStoreUnionMembersAsync" }       # SUNIONSTORE
    @{ Row = 167; StringMethod = $null;                                   Method = "SetIsMemberAsync" }                # SISMEMBER
    @{ Row = 168; StringMethod = "SetGetMembersStringAsync";              Method = "SetGetMembersAsync" }              # SMEMBERS
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $true
    if ($item.StringMethod) {
        $ws.Cells.Item($r, 4).Value = $item.StringMethod
    }
    $ws.Cells.Item($r, 5).Value = $item.Method
}

# Scroll position / final selection left by the editor after the change.
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 144
    $win.ScrollColumn = 1
} catch {
    # view scroll position is cosmetic only; ignore if unsupported
}
$ws.Range("C173:C174").Select()
